$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Range("A44").Value = 111901529
$ws.Range("B44").Value = 56398
$ws.Range("D44").Value = 'NT'
$ws.Range("E44").Value = 100109
$ws.Range("F44").Value = 'Tretåig hackspett'
$ws.Range("G44").Value = 'Picoides tridactylus'
$ws.Range("H44").Value = '(Linnaeus, 1758)'
$ws.Range("Q44").Value = 478295.8274075754
$ws.Range("R44").Value = 7034510.601185531
$ws.Range("AC44").Value = 'ringhack färska'

# Row 45
$ws.Range("A45").Value = 111901607
$ws.Range("B45").Value = 89423
$ws.Range("D45").Value = 'NT'
$ws.Range("E45").Value = 5432
$ws.Range("F45").Value = 'Granticka'
$ws.Range("G45").Value = 'Porodaedalea chrysoloma'
$ws.Range("H45").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q45").Value = 478095.1015727581
$ws.Range("R45").Value = 7035043.219008418
$ws.Range("AC45").Value = ""

# Row 46
$ws.Range("A46").Value = 111901554
$ws.Range("B46").Value = 56398
$ws.Range("D46").Value = 'NT'
$ws.Range("E46").Value = 100109
$ws.Range("F46").Value = 'Tretåig hackspett'
$ws.Range("G46").Value = 'Picoides tridactylus'
$ws.Range("H46").Value = '(Linnaeus, 1758)'
$ws.Range("Q46").Value = 477059.476171807
$ws.Range("R46").Value = 7033542.438482954
$ws.Range("AC46").Value = 'ringhack äldre'

# Row 47
$ws.Range("A47").Value = 111901578
$ws.Range("B47").Value = 56398
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 100109
$ws.Range("F47").Value = 'Tretåig hackspett'
$ws.Range("G47").Value = 'Picoides tridactylus'
$ws.Range("H47").Value = '(Linnaeus, 1758)'
$ws.Range("Q47").Value = 477843.0506277476
$ws.Range("R47").Value = 7034173.07203023
$ws.Range("AC47").Value = 'ringhack'

# Row 48
$ws.Range("A48").Value = 111901574
$ws.Range("B48").Value = 56398
$ws.Range("D48").Value = 'NT'
$ws.Range("E48").Value = 100109
$ws.Range("F48").Value = 'Tretåig hackspett'
$ws.Range("G48").Value = 'Picoides tridactylus'
$ws.Range("H48").Value = '(Linnaeus, 1758)'
$ws.Range("Q48").Value = 477521.0595750482
$ws.Range("R48").Value = 7034024.014202636
$ws.Range("AC48").Value = 'ringhack äldre'

# Row 49
$ws.Range("A49").Value = 111901622
$ws.Range("B49").Value = 85062
$ws.Range("D49").Value = 'NT'
$ws.Range("E49").Value = 249278
$ws.Range("F49").Value = 'Barrviolspindling'
$ws.Range("G49").Value = 'Cortinarius harcynicus'
$ws.Range("H49").Value = '(Pers.) M.M.Moser'
$ws.Range("Q49").Value = 478165.647914707
$ws.Range("R49").Value = 7034284.10291774
$ws.Range("AC49").Value = ""

# Row 50
$ws.Range("A50").Value = 111901571
$ws.Range("B50").Value = 56398
$ws.Range("D50").Value = 'NT'
$ws.Range("E50").Value = 100109
$ws.Range("F50").Value = 'Tretåig hackspett'
$ws.Range("G50").Value = 'Picoides tridactylus'
$ws.Range("H50").Value = '(Linnaeus, 1758)'
$ws.Range("Q50").Value = 477388.9837839347
$ws.Range("R50").Value = 7033793.496102724
$ws.Range("AC50").Value = 'ringhack färska'

# Row 54
$ws.Range("A54").Value = 111901564
$ws.Range("B54").Value = 56398
$ws.Range("D54").Value = 'NT'
$ws.Range("E54").Value = 100109
$ws.Range("F54").Value = 'Tretåig hackspett'
$ws.Range("G54").Value = 'Picoides tridactylus'
$ws.Range("H54").Value = '(Linnaeus, 1758)'
$ws.Range("Q54").Value = 477365.3609946552
$ws.Range("R54").Value = 7033686.214811271
$ws.Range("AC54").Value = 'ringhack'

# Row 55
$ws.Range("A55").Value = 111901609
$ws.Range("B55").Value = 77515
$ws.Range("D55").Value = 'NT'
$ws.Range("E55").Value = 6425
$ws.Range("F55").Value = 'Garnlav'
$ws.Range("G55").Value = 'Alectoria sarmentosa'
$ws.Range("H55").Value = '(Ach.) Ach.'
$ws.Range("Q55").Value = 478140.6856806503
$ws.Range("R55").Value = 7034828.538396582
$ws.Range("AC55").Value = ""

# Row 56
$ws.Range("A56").Value = 111901573
$ws.Range("B56").Value = 56398
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 100109
$ws.Range("F56").Value = 'Tretåig hackspett'
$ws.Range("G56").Value = 'Picoides tridactylus'
$ws.Range("H56").Value = '(Linnaeus, 1758)'
$ws.Range("Q56").Value = 477537.5564934253
$ws.Range("R56").Value = 7034011.363671634
$ws.Range("AC56").Value = 'ringhack färska'

# Row 57
$ws.Range("A57").Value = 111901561
$ws.Range("B57").Value = 56398
$ws.Range("D57").Value = 'NT'
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("Q57").Value = 477353.6324963307
$ws.Range("R57").Value = 7033676.000540117
$ws.Range("AC57").Value = 'ringhack'

# Row 58
$ws.Range("A58").Value = 111901532
$ws.Range("B58").Value = 56398
$ws.Range("D58").Value = 'NT'
$ws.Range("E58").Value = 100109
$ws.Range("F58").Value = 'Tretåig hackspett'
$ws.Range("G58").Value = 'Picoides tridactylus'
$ws.Range("H58").Value = '(Linnaeus, 1758)'
$ws.Range("Q58").Value = 478222.3906325128
$ws.Range("R58").Value = 7034454.703636711
$ws.Range("AC58").Value = 'ringhack'

# Row 59
$ws.Range("A59").Value = 111901534
$ws.Range("B59").Value = 56398
$ws.Range("D59").Value = 'NT'
$ws.Range("E59").Value = 100109
$ws.Range("F59").Value = 'Tretåig hackspett'
$ws.Range("G59").Value = 'Picoides tridactylus'
$ws.Range("H59").Value = '(Linnaeus, 1758)'
$ws.Range("Q59").Value = 478196.6579575058
$ws.Range("R59").Value = 7034427.575356619
$ws.Range("AC59").Value = 'ringhack'

# Row 60
$ws.Range("A60").Value = 111901523
$ws.Range("B60").Value = 90087
$ws.Range("D60").Value = 'LC'
$ws.Range("E60").Value = 3298
$ws.Range("F60").Value = 'Trådticka'
$ws.Range("G60").Value = 'Climacocystis borealis'
$ws.Range("H60").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q60").Value = 478095.1199801999
$ws.Range("R60").Value = 7035045.903991099
$ws.Range("AC60").Value = ""

# Row 78
$ws.Range("A78").Value = 111901555
$ws.Range("B78").Value = 56398
$ws.Range("D78").Value = 'NT'
$ws.Range("E78").Value = 100109
$ws.Range("F78").Value = 'Tretåig hackspett'
$ws.Range("G78").Value = 'Picoides tridactylus'
$ws.Range("H78").Value = '(Linnaeus, 1758)'
$ws.Range("Q78").Value = 477075.053782316
$ws.Range("R78").Value = 7033651.571049522
$ws.Range("AC78").Value = 'ringhack'

# Row 79
$ws.Range("A79").Value = 111901590
$ws.Range("B79").Value = 56414
$ws.Range("D79").Value = 'NT'
$ws.Range("E79").Value = 100049
$ws.Range("F79").Value = 'Spillkråka'
$ws.Range("G79").Value = 'Dryocopus martius'
$ws.Range("H79").Value = '(Linnaeus, 1758)'
$ws.Range("Q79").Value = 477995.937568082
$ws.Range("R79").Value = 7034178.282797099
$ws.Range("AC79").Value = 'hack'

# Row 80
$ws.Range("A80").Value = 111901589
$ws.Range("B80").Value = 90092
$ws.Range("D80").Value = 'VU'
$ws.Range("E80").Value = 67
$ws.Range("F80").Value = 'Sprickporing'
$ws.Range("G80").Value = 'Diplomitoporus crustulinus'
$ws.Range("H80").Value = '(Bres.) Domański'
$ws.Range("Q80").Value = 477449.9680636173
$ws.Range("R80").Value = 7033732.178319109
$ws.Range("AC80").Value = ""

# Row 81
$ws.Range("A81").Value = 111901530
$ws.Range("B81").Value = 56398
$ws.Range("D81").Value = 'NT'
$ws.Range("E81").Value = 100109
$ws.Range("F81").Value = 'Tretåig hackspett'
$ws.Range("G81").Value = 'Picoides tridactylus'
$ws.Range("H81").Value = '(Linnaeus, 1758)'
$ws.Range("Q81").Value = 478301.071792486
$ws.Range("R81").Value = 7034490.871451757
$ws.Range("AC81").Value = 'ringhack äldre'
